# Minor issues caught in quality checks
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 ("Know your tools: Example"): clarify caption about library paths
# and widen the auto-fit textbox so the longer caption still fits on one line.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$caption = $s6.Shapes.Item(15)
$caption.TextFrame.TextRange.Text = "Yellow ellipses in library paths indicate snipped content"
$caption.Width = 361.6829133858268

# ---------------------------------------------------------------------------
# Slide 18 ("Is this working?"): label the closing hyperlink and tidy up the
# spacer paragraph above it (turn off its bullet now that it is blank).
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$body18 = $s18.Shapes.Item(2)
$tr18 = $body18.TextFrame.TextRange

$spacerPara = $tr18.Paragraphs(8)
$spacerPara.ParagraphFormat.Bullet.Type = 0

$linkPara = $tr18.Paragraphs(9)
$linkPara.InsertBefore("For your exploration: ")
$newLead = $tr18.Characters($linkPara.Start, 22)
$newLead.ActionSettings.Item(1).Hyperlink.Address = ""

# ---------------------------------------------------------------------------
# Slide 24 ("How do you plan"): fix capitalization of "vs".
# ---------------------------------------------------------------------------
$s24 = $p.Slides.Item(24)
$body24 = $s24.Shapes.Item(2)
$tr24 = $body24.TextFrame.TextRange
$costPara = $tr24.Paragraphs(6)
$costPara.Text = "Cost benefit analysis of fidelity vs reaching science goals in allocated resources"
